$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set a cell to a string value, forcing text storage even when the string
# looks like a plain number (e.g. "76.09"), while leaving the cell's style
# back at the default ("Normal") once the text is committed so no stray
# number-format style lingers on the cell.
function Set-TextCell {
    param(
        [string]$Addr,
        [string]$Val
    )
    $cell = $ws.Range($Addr)
    $cell.NumberFormat = "@"
    $cell.Value = $Val
    $cell.Style = "Normal"
}

function Set-Row {
    param(
        [int]$Row,
        [string]$B,
        [string]$C,
        [string]$D,
        [string]$E
    )
    if ($B -ne "") { $ws.Range("B$Row").Value = $B }
    if ($C -ne "") { $ws.Range("C$Row").Value = $C }
    if ($D -ne "") { Set-TextCell "D$Row" $D }
    if ($E -ne "") { $ws.Range("E$Row").Value = $E }
}

Set-Row 2  "" "" "42.403.84"  "  +0.34%  "
Set-Row 3  "" "" "2.245.64"   "  -0.03%  "
Set-Row 4  "" "" ""           "  -0.03%  "
Set-Row 5  "" "" "246.53"     "  -0.31%  "
Set-Row 6  "" "" "0.624"      "  -1.22%  "
Set-Row 7  "" "" "76.09"      "  -0.17%  "
Set-Row 8  "" "" ""           "  +0.01%  "
Set-Row 9  "" "" ""           "  -1.38%  "
Set-Row 10 "" "" "44.10"      "  +10.00%  "
Set-Row 11 "" "" "0.0944"     "  -0.69%  "
Set-Row 12 "" "" "7.32"       "  +1.84%  "
Set-Row 13 "" "" ""           "  -1.08%  "
Set-Row 14 "" "" "2.588.35"   "  +0.21%  "
Set-Row 15 "" "" "14.62"      "  -1.94%  "
Set-Row 17 "" "" "2.235.46"   "  -1.18%  "
Set-Row 18 "" "" "42.217.70"  "  -0.03%  "
Set-Row 19 "" "" ""           "  +4.14%  "
Set-Row 20 "" "" ""           "  +0.59%  "
Set-Row 21 "" "" "72.24"      "  +0.90%  "
Set-Row 22 "" "" ""           "  +2.66%  "
Set-Row 23 "" "" "231.63"     "  +0.04%  "
Set-Row 24 "" "" "9.15"       "  +31.22%  "
Set-Row 25 "" "" ""           "  +0.04%  "
Set-Row 26 "" "" "11.46"      "  +3.41%  "
Set-Row 27 "" "" ""           "  -3.16%  "
Set-Row 28 "" "" ""           "  -0.52%  "
Set-Row 29 "" "" "2.20"       "  +1.38%  "
Set-Row 30 "" "" "168.19"     "  +0.13%  "
Set-Row 31 "" "" "20.70"      "  +0.80%  "
Set-Row 32 "" "" "0.0832"     ""
Set-Row 33 "" "" "0.120"      "  +0.53%  "
Set-Row 34 "" "" "30.64"      "  -3.76%  "

# Rows 35 and 36 swap content (Stellar <-> Filecoin)
Set-Row 35 "Filecoin" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil" "5.29" "  +11.77%  "
Set-Row 36 "Stellar"  "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm" "0.125" "  -0.52%  "

Set-Row 37 "" "" "4.55"      "  +1.42%  "
Set-Row 38 "" "" "0.0319"    "  +7.21%  "
Set-Row 39 "" "" "13.72"     "  +5.83%  "
Set-Row 40 "" "" ""          "  -1.81%  "
Set-Row 41 "" "" "5.81"      "  -1.88%  "
Set-Row 42 "" "" "63.69"     "  +5.88%  "
Set-Row 43 "" "" ""          "  -0.45%  "
Set-Row 44 "" "" "108.12"    "  -8.22%  "
Set-Row 45 "" "" ""          "  +0.64%  "
Set-Row 46 "" "" ""          "  +1.49%  "
Set-Row 47 "" "" "0.997"     "  -0.14%  "
Set-Row 48 "" "" ""          "  +0.71%  "
Set-Row 49 "" "" ""          "  -0.47%  "
Set-Row 50 "" "" ""          "  +6.16%  "

Set-Row 51 "WOONetwork" "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo" "0.429" "  +19.49%  "
